# Weekly update: insert a new week's price record for Choclo (Vega Monumental
# Concepción) as row 30, shifting every subsequent record down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 30:50 down to 31:51 by inserting a new blank row at 30.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly data point.
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44447
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = 100112024
$ws.Range("G30").Value = "Choclo"
$ws.Range("H30").Value = "Dulce o Americano"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 35000
$ws.Range("L30").Value = 36000
$ws.Range("M30").Value = 35500
$ws.Range("N30").Value = "$/malla 70 unidades"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 507
$ws.Range("Q30").Value = 70
$ws.Range("R30").Value = "Hortaliza"
